$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly Qatar M2 data points appended to the bottom of the table.
# Row 196
$ws.Range("A195").Copy()
$ws.Range("A196").PasteSpecial(-4122)
$ws.Range("A196").Value = 44986.45833333334
$ws.Range("B196").Value = "ECONOMICS:QAM2"
$ws.Range("C196").Value = 704618000000
$ws.Range("D196").Value = 704618000000
$ws.Range("E196").Value = 704618000000
$ws.Range("F196").Value = 704618000000
$ws.Range("G196").Value = 0

# Row 197
$ws.Range("A195").Copy()
$ws.Range("A197").PasteSpecial(-4122)
$ws.Range("A197").Value = 45017.45833333334
$ws.Range("B197").Value = "ECONOMICS:QAM2"
$ws.Range("C197").Value = 696270000000
$ws.Range("D197").Value = 696270000000
$ws.Range("E197").Value = 696270000000
$ws.Range("F197").Value = 696270000000
$ws.Range("G197").Value = 0

$excel.CutCopyMode = 0
